$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("LandByRegion")
$ws3 = $wb.Worksheets.Item("A of MEA")

# ---------------------------------------------------------------------------
# Append a duplicate of the "A of MEA" table (A1:E12) onto the "LandByRegion"
# sheet at rows 16-27 (mirrors copy/pasting that block of data further down
# the first worksheet).
# ---------------------------------------------------------------------------

# Row 16 - header row
$ws1.Range("A16").Value = "Country Name"
$ws1.Range("B16").NumberFormat = "@"
$ws1.Range("B16").Value = "Country Code"
$ws1.Range("C16").Value = "Series Name"
$ws1.Range("D16").NumberFormat = "@"
$ws1.Range("D16").Value = "Series Code"
$ws1.Range("E16").Value = "2018 [YR2018]"

# Row 17 - Algeria
$ws1.Range("A17").Value = "Algeria"
$ws1.Range("B17").NumberFormat = "@"
$ws1.Range("B17").Value = "DZA"
$ws1.Range("C17").Value = "Land area (sq. km)"
$ws1.Range("D17").NumberFormat = "@"
$ws1.Range("D17").Value = "AG.LND.TOTL.K2"
$ws1.Range("E17").Value = 2381740

# Row 18 - Egypt, Arab Rep.
$ws1.Range("A18").Value = "Egypt, Arab Rep."
$ws1.Range("B18").NumberFormat = "@"
$ws1.Range("B18").Value = "EGY"
$ws1.Range("C18").Value = "Land area (sq. km)"
$ws1.Range("D18").NumberFormat = "@"
$ws1.Range("D18").Value = "AG.LND.TOTL.K2"
$ws1.Range("E18").Value = 995450

# Row 19 - Morocco
$ws1.Range("A19").Value = "Morocco"
$ws1.Range("B19").NumberFormat = "@"
$ws1.Range("B19").Value = "MAR"
$ws1.Range("C19").Value = "Land area (sq. km)"
$ws1.Range("D19").NumberFormat = "@"
$ws1.Range("D19").Value = "AG.LND.TOTL.K2"
$ws1.Range("E19").Value = 446300

# Row 20 - Libya
$ws1.Range("A20").Value = "Libya"
$ws1.Range("B20").NumberFormat = "@"
$ws1.Range("B20").Value = "LBY"
$ws1.Range("C20").Value = "Land area (sq. km)"
$ws1.Range("D20").NumberFormat = "@"
$ws1.Range("D20").Value = "AG.LND.TOTL.K2"
$ws1.Range("E20").Value = 1759540

# Row 21 - Djibouti
$ws1.Range("A21").Value = "Djibouti"
$ws1.Range("B21").NumberFormat = "@"
$ws1.Range("B21").Value = "DJI"
$ws1.Range("C21").Value = "Land area (sq. km)"
$ws1.Range("D21").NumberFormat = "@"
$ws1.Range("D21").Value = "AG.LND.TOTL.K2"
$ws1.Range("E21").Value = 23180

# Row 22 - Tunisia
$ws1.Range("A22").Value = "Tunisia"
$ws1.Range("B22").NumberFormat = "@"
$ws1.Range("B22").Value = "TUN"
$ws1.Range("C22").Value = "Land area (sq. km)"
$ws1.Range("D22").NumberFormat = "@"
$ws1.Range("D22").Value = "AG.LND.TOTL.K2"
$ws1.Range("E22").Value = 155360

# Row 23 - totals row (formula)
$ws1.Range("B23").NumberFormat = "@"
$ws1.Range("D23").NumberFormat = "@"
$ws1.Range("E23").Formula = "=SUM(E17:E22)"

# Rows 24 & 25 - blank spacer rows (keep the Text format in B/D like the source)
$ws1.Range("B24").NumberFormat = "@"
$ws1.Range("D24").NumberFormat = "@"
$ws1.Range("B25").NumberFormat = "@"
$ws1.Range("D25").NumberFormat = "@"

# Row 26 - source note
$ws1.Range("A26").Value = "Data from database: World Development Indicators"
$ws1.Range("B26").NumberFormat = "@"
$ws1.Range("D26").NumberFormat = "@"

# Row 27 - last-updated note
$ws1.Range("A27").Value = "Last Updated: 07/01/2020"

# ---------------------------------------------------------------------------
# View/selection state
# ---------------------------------------------------------------------------

# "A of MEA" ends up with A1:F13 selected (but stays inactive - LandByRegion
# remains the visible tab).
$ws3.Range("A1:F13").Select()
$ws1.Activate()

# LandByRegion itself ends up scrolled down with L21 as the active cell.
$ws1.Range("L21").Select()

# ---------------------------------------------------------------------------
# Page setup on "A of MEA"
# ---------------------------------------------------------------------------
$ws3.PageSetup.PaperSize = 9
$ws3.PageSetup.Orientation = 1
